$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MapNpcInfo")
$ws2 = $wb.Worksheets.Item("MapNpcPosition")
$ws3 = $wb.Worksheets.Item("MapNpcMenu")

# --- MapNpcInfo: insert a new column C (ResourceKey / string) ---
$ws1.Columns.Item(3).Insert()
$ws1.Columns.Item(3).ColumnWidth = 24.14

$c1 = $ws1.Range("C1")
$c1.Font.Name = "Roboto"
$c1.Font.Size = 10
$c1.Font.Color = 0
$c1.Interior.Color = 15724527
$c1.Borders.LineStyle = 0
$c1.HorizontalAlignment = -4131
$c1.Value = "string"

$c2 = $ws1.Range("C2")
$c2.Font.Name = "Roboto"
$c2.Font.Size = 10
$c2.Font.Color = 0
$c2.Interior.Color = 15724527
$c2.Borders.LineStyle = 0
$c2.Value = "ResourceKey"

$ws1.Range("C4").Value = "img_powder_shop_mirror"

# Myhome slime sprite position/scale value changed
$ws1.Range("D3").Value = "0,36,0.5"

# --- selections / active sheet ---
$ws1.Range("E2").Select()
$ws3.Range("C14").Select()
$ws2.Activate()
